$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the data rows: "Dongle" master data records replaced with "Desktop Computer" records.
$ws.Range("A2").Value = "DKS"
$ws.Range("B2").Value = "Dekstop"
$ws.Range("C2").Value = "Desktop Computer"

$ws.Range("A3").Value = "DKS"
$ws.Range("B3").Value = "الحاسوب"
$ws.Range("C3").Value = "أجهزة الكمبيوتر المكتبية"

$ws.Range("A4").Value = "DKS"
$ws.Range("B4").Value = "Ordinateur"
$ws.Range("C4").Value = "Ordinateurs de bureau"

# Record the active cell/selection as it was when the file was saved.
$ws.Range("D10").Select()

# Configure page setup (adds <pageSetup .../> to the worksheet XML).
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

$wb.Save()
